$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    ,@("Basica", "2022-09-14", 428)
    ,@("Media", "2022-09-14", 576)
    ,@("UTI", "2022-09-14", 208)
    ,@("UCI", "2022-09-14", 131)
    ,@("Basica", "2022-09-15", 428)
    ,@("Media", "2022-09-15", 555)
    ,@("UTI", "2022-09-15", 203)
    ,@("UCI", "2022-09-15", 136)
    ,@("Basica", "2022-09-16", 431)
    ,@("Media", "2022-09-16", 524)
    ,@("UTI", "2022-09-16", 203)
    ,@("UCI", "2022-09-16", 130)
    ,@("Basica", "2022-09-17", 405)
    ,@("Media", "2022-09-17", 493)
    ,@("UTI", "2022-09-17", 194)
    ,@("UCI", "2022-09-17", 131)
    ,@("Basica", "2022-09-18", 410)
    ,@("Media", "2022-09-18", 496)
    ,@("UTI", "2022-09-18", 202)
    ,@("UCI", "2022-09-18", 131)
    ,@("Basica", "2022-09-19", 417)
    ,@("Media", "2022-09-19", 512)
    ,@("UTI", "2022-09-19", 179)
    ,@("UCI", "2022-09-19", 131)
    ,@("Basica", "2022-09-20", 401)
    ,@("Media", "2022-09-20", 528)
    ,@("UTI", "2022-09-20", 193)
    ,@("UCI", "2022-09-20", 135)
    ,@("Basica", "2022-09-21", 427)
    ,@("Media", "2022-09-21", 546)
    ,@("UTI", "2022-09-21", 191)
    ,@("UCI", "2022-09-21", 139)
    ,@("Basica", "2022-09-22", 421)
    ,@("Media", "2022-09-22", 542)
    ,@("UTI", "2022-09-22", 187)
    ,@("UCI", "2022-09-22", 138)
)

$startRow = 3526
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $item = $data[$i]
    $ws.Cells.Item($row, 1).Value = $item[0]
    $dateCell = $ws.Cells.Item($row, 2)
    $dateCell.NumberFormat = "@"
    $dateCell.Value = $item[1]
    $dateCell.Style = "Normal"
    $ws.Cells.Item($row, 3).Value = $item[2]
}
